$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Query text used in both cells: replace "IN ['HISPANIC_OR_LATINO']" with '= "HISPANIC_OR_LATINO"'
$b2Text = @"
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity = "HISPANIC_OR_LATINO"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``
"@

$c2Text = @"
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE WHERE c.ethnicity = "HISPANIC_OR_LATINO"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
"@

$ws.Range("B2").Value = $b2Text
$ws.Range("C2").Value = $c2Text

# New row 3: empty cells carrying the same wrap-text style as B2/C2
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Update selection to C2
$ws.Range("C2").Select()
